$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill row 7 with the new 2050 Export Demand entry, mirroring rows 5/6.
$ws.Range("D7").Value = 2050
$ws.Range("E7").Value = "LO"
$ws.Range("F7").Value = "ACT_BND"
$ws.Range("H7").Value = 50
$ws.Range("J7").Value = 0
$ws.Range("M7").Value = "EXPH2*"

# Apply the same number-format styling as the rows above (G:I use style 14 /
# "Comma" numFmt 43, J:L use style 15 / same numFmt with a border).
$ws.Range("G7:I7").Style = $ws.Range("G6:I6").Style
$ws.Range("J7:L7").Style = $ws.Range("J6:L6").Style
$ws.Range("M7").Style = $ws.Range("M6").Style

# Move the active selection to S14, as recorded in the workbook view.
$ws.Range("S14").Select()
